$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new text values look numeric, so Excel keeps them as text
# (matches original workbook where these Price cells are stored as text, not numbers)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values from the crypto data refresh
$ws.Range("D2").Value = "28.613.36"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.583.79"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "213.63"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "44.24"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "24.13"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.810.74"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "1.581.12"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "28.623.24"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "62.28"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "231.87"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "2.07"
$ws.Range("E25").Value = "  +4.76%  "
$ws.Range("D26").Value = "151.73"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "15.07"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "1.400.64"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "2.65"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "0.523"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "1.90"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("D47").Value = "0.962"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").Value = "63.34"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "1.722.29"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "86.87"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  -2.20%  "
